# Swap the deck's applied theme (ppt/theme/theme1.xml, currently the
# "Integral" / "Red Violet" color scheme used by the slide master) over to
# the standard Office theme colors, matching the target commit which
# exchanged the contents of theme1.xml and theme2.xml (font scheme and
# format scheme are already identical between the two theme parts, so the
# only meaningful, renderable difference is the 12-slot color scheme).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# PpThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1-6, 11 hlink, 12 folHlink.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
